$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Rename Sheet1 -> Lenovo
$ws1.Name = "Lenovo"

# Set the print area on the renamed sheet (adds the
# _xlnm.Print_Area defined name scoped to this sheet)
$ws1.PageSetup.PrintArea = "A1:L28"

# Update the three embedded charts so their series formulas
# point at the renamed sheet instead of the old "Sheet1" name
$chartObjects = $ws1.ChartObjects()

$chart1 = $chartObjects.Item(1).Chart
$series1 = $chart1.SeriesCollection(1)
$series1.XValues = "=Lenovo!`$A`$2:`$A`$21"
$series1.Values = "=Lenovo!`$H`$2:`$H`$21"

$chart2 = $chartObjects.Item(2).Chart
$series2 = $chart2.SeriesCollection(1)
$series2.XValues = "=Lenovo!`$A`$2:`$A`$21"
$series2.Values = "=Lenovo!`$K`$2:`$K`$21"

$chart3 = $chartObjects.Item(3).Chart
$series3 = $chart3.SeriesCollection(1)
$series3.Name = "=Lenovo!`$L`$1"
$series3.XValues = "=Lenovo!`$A`$2:`$A`$21"
$series3.Values = "=Lenovo!`$L`$2:`$L`$21"

# Grow the header row
$ws1.Rows.Item(1).RowHeight = 90

# Match the new selection left behind after defining the print area
$null = $ws1.Range("A1:L28").Select()
